# Add season record columns (Wins, Losses, Ties) to the team stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy the existing header formatting (bold + border, style s="1")
# onto the three new header cells, then set their labels.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-50: every player row on this roster shares the same team
# season record, so fill the same Wins/Losses/Ties values down the column.
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = 79   # AD -> Wins
    $ws.Cells.Item($row, 31).Value = 83   # AE -> Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF -> Ties
}
